$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.768.89"
$ws.Range("E2").Value = "  -1.90%  "
$ws.Range("D3").Value = "1.610.38"
$ws.Range("E3").Value = "  -3.95%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.81"
$ws.Range("E5").Value = "  -2.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5192"
$ws.Range("E6").Value = "  -1.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2550"
$ws.Range("E9").Value = "  -1.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.05"
$ws.Range("E10").Value = "  -6.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07514"
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("D12").Value = "1.620.87"
$ws.Range("E12").Value = "  -3.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.338"
$ws.Range("E13").Value = "  -3.04%  "
$ws.Range("D14").Value = "1.839.95"
$ws.Range("E14").Value = "  -3.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5384"
$ws.Range("E15").Value = "  -4.51%  "
$ws.Range("D16").Value = "0.0₅7798"
$ws.Range("E16").Value = "  -2.94%  "
$ws.Range("E17").Value = "  -4.99%  "
$ws.Range("D18").Value = "25.775.61"
$ws.Range("E18").Value = "  -1.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.004"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.574"
$ws.Range("E20").Value = "  -5.37%  "
$ws.Range("E21").Value = "  -3.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.954"
$ws.Range("E22").Value = "  -4.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.004"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("E24").Value = "  -3.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.24"
$ws.Range("E25").Value = "  -3.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1198"
$ws.Range("E26").Value = "  -4.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.298"
$ws.Range("E27").Value = "  -3.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.41"
$ws.Range("E28").Value = "  -3.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.357"
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05885"
$ws.Range("E30").Value = "  -5.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.232"
$ws.Range("E31").Value = "  -4.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.355"
$ws.Range("E32").Value = "  -4.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.303"
$ws.Range("E33").Value = "  -4.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.581"
$ws.Range("E34").Value = "  -3.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9580"
$ws.Range("E35").Value = "  -4.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.381"
$ws.Range("E36").Value = "  -1.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.694"
$ws.Range("E37").Value = "  -1.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5656"
$ws.Range("E38").Value = "  -6.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01576"
$ws.Range("E39").Value = "  -2.94%  "
$ws.Range("E40").Value = "  -0.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8310"
$ws.Range("E41").Value = "  -4.65%  "
$ws.Range("D42").Value = "1.015.73"
$ws.Range("E42").Value = "  -8.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.582"
$ws.Range("E43").Value = "  -8.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "98.78"
$ws.Range("E44").Value = "  -1.16%  "
$ws.Range("D45").Value = "1.762.47"
$ws.Range("E45").Value = "  -3.54%  "
$ws.Range("E46").Value = "  -1.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9997"
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.78"
$ws.Range("E48").Value = "  -4.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05154"
$ws.Range("E49").Value = "  -1.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.810"
$ws.Range("E50").Value = "  -2.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4210"
$ws.Range("E51").Value = "  -1.08%  "
